$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Gracie Murphy's bio (row 6, column H)
$ws.Range("H6").Value = 'I am currently a 3rd year an MIMG major with a minor in Professional Writing. I have been working with the Jacob’s Laboratory for the past year researching the impact of the human gut microbiome on stress and human disease. I love working with kids and am excited to work with local elementary and middle schools this upcoming year.'

# Update Tyler Wu's bio (row 9, column H)
$ws.Range("H9").Value = 'Hey y''all. I''m a fourth year data theory and cognitive science major. In my free time I enjoy, drawing, coding, and plyaing basketball. I''m always interested in learning more about a variety of subjects, including biology. Thus, I hope this year I can learn more about virology, while contributing to the creative side of this organization.'

# Update the view/selection state to match the new event format
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("I7").Select()
